$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.066.95"
$ws.Range("E2").Value = "  +2.59%  "

$ws.Range("D3").Value = "2.304.29"
$ws.Range("E3").Value = "  +2.10%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'302.39"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("D6").Value = "'99.42"

$ws.Range("E7").Value = "  +1.92%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +3.28%  "

$ws.Range("E10").Value = "  +4.66%  "

$ws.Range("E11").Value = "  +1.71%  "

$ws.Range("D12").Value = "'49.07"
$ws.Range("E12").Value = "  +3.41%  "

$ws.Range("E13").Value = "  +4.28%  "

$ws.Range("D14").Value = "'17.79"
$ws.Range("E14").Value = "  +16.87%  "

$ws.Range("E15").Value = "  +2.42%  "

$ws.Range("D16").Value = "2.663.36"
$ws.Range("E16").Value = "  +2.13%  "

$ws.Range("D17").Value = "2.292.68"
$ws.Range("E17").Value = "  +1.57%  "

$ws.Range("E18").Value = "  +4.58%  "

$ws.Range("D19").Value = "42.956.57"
$ws.Range("E19").Value = "  +2.27%  "

$ws.Range("D20").Value = "'12.36"
$ws.Range("E20").Value = "  +9.09%  "

$ws.Range("E21").Value = "  +2.20%  "

$ws.Range("E22").Value = "  +1.68%  "

$ws.Range("D23").Value = "'67.82"
$ws.Range("E23").Value = "  +2.00%  "

$ws.Range("D24").Value = "'237.06"
$ws.Range("E24").Value = "  +1.84%  "

$ws.Range("D25").Value = "'2.17"
$ws.Range("E25").Value = "  +13.41%  "

$ws.Range("E26").Value = "  +1.19%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").Value = "'24.56"
$ws.Range("E28").Value = "  +3.68%  "

$ws.Range("D29").Value = "'168.32"
$ws.Range("E29").Value = "  +1.15%  "

$ws.Range("E30").Value = "  -3.03%  "

$ws.Range("D31").Value = "'33.87"
$ws.Range("E31").Value = "  +1.49%  "

$ws.Range("D32").Value = "'9.18"
$ws.Range("E32").Value = "  +1.83%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("D34").Value = "'5.01"
$ws.Range("E34").Value = "  +1.87%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "'4.56"
$ws.Range("E35").Value = "  +4.89%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.43"
$ws.Range("E36").Value = "  +3.71%  "

$ws.Range("D37").Value = "'17.04"
$ws.Range("E37").Value = "  +7.65%  "

$ws.Range("D38").Value = "'0.0701"
$ws.Range("E38").Value = "  +1.35%  "

$ws.Range("E39").Value = "  +3.78%  "

$ws.Range("E40").Value = "  +5.16%  "

$ws.Range("E41").Value = "  +1.06%  "

$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("E43").Value = "  -2.26%  "

$ws.Range("D44").Value = "2.002.76"
$ws.Range("E44").Value = "  +3.07%  "

$ws.Range("D45").Value = "'0.0285"
$ws.Range("E45").Value = "  +2.91%  "

$ws.Range("D46").Value = "'10.12"
$ws.Range("E46").Value = "  +6.18%  "

$ws.Range("D47").Value = "'17.68"
$ws.Range("E47").Value = "  +2.69%  "

$ws.Range("E48").Value = "  +3.29%  "

$ws.Range("D49").Value = "'55.49"
$ws.Range("E49").Value = "  +6.66%  "

$ws.Range("D50").Value = "2.529.45"
$ws.Range("E50").Value = "  +1.89%  "

$ws.Range("E51").Value = "  +2.69%  "
